$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Header cell: "Characteristic" -> "Baseline Characteristics", no longer bold ---
$headerCell = $t.Rows.Item(1).Cells.Item(1)
$headerCell.Range.Text = "Baseline Characteristics"

# Restrict the bold change to the text run itself (exclude the trailing paragraph
# mark) so only the run's own formatting is touched, not the paragraph mark's.
$headerRange = $headerCell.Range
$headerTextRange = $d.Range($headerRange.Start, $headerRange.End - 1)
$headerTextRange.Font.Bold = $false

# --- 2/3. Swap the "Greenspace 300m ..." and "Water 300m ..." rows (text + row height) ---
# Row 34 currently holds the Greenspace text/height, Row 35 the Water text/height;
# the edit swaps the two rows' label text and trHeight while leaving the numeric
# value cells untouched.
$greenRow = $t.Rows.Item(34)
$waterRow = $t.Rows.Item(35)

# Stage through a unique placeholder so the two cells never briefly hold identical
# text at the same time (avoids the two ranges getting conflated mid-swap).
$waterRow.Cells.Item(1).Range.Text = "__SWAP_PLACEHOLDER__"
$greenRow.Cells.Item(1).Range.Text = "Water 300m from residence, Median (Q1, Q3)"
$waterRow.Cells.Item(1).Range.Text = "Greenspace 300m from residence, Median (Q1, Q3)"

$greenRow.Height = 30.35
$waterRow.Height = 31.2
